# Auto-generated edit script applying numeric updates described in the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 31251274
$ws.Range("I6").Value = 999
$ws.Range("J6").Value = 41668030
$ws.Range("K6").Value = 2997
$ws.Range("L6").Value = 125004090
$ws.Range("M6").Value = -2885
$ws.Range("N6").Value = -125004314
$ws.Range("H17").Value = 20001998
$ws.Range("J17").Value = 20001998
$ws.Range("L17").Value = 60005994
$ws.Range("N17").Value = -60006330
$ws.Range("H28").Value = 217.11111
$ws.Range("I28").Value = 224
$ws.Range("K28").Value = 224
$ws.Range("M28").Value = 261
$ws.Range("H33").Value = 663.0714
$ws.Range("I33").Value = 607.6667
$ws.Range("J33").Value = 995.5
$ws.Range("K33").Value = 607.6667
$ws.Range("L33").Value = 995.5
$ws.Range("M33").Value = -378.6667
$ws.Range("N33").Value = -1453.5
$ws.Range("H40").Value = 4778.8887
$ws.Range("J40").Value = 5267.857
$ws.Range("L40").Value = 5267.857
$ws.Range("N40").Value = -5617.857
$ws.Range("H64").Value = 11834
$ws.Range("I64").Value = 6446.6665
$ws.Range("K64").Value = 6446.6665
$ws.Range("M64").Value = -6198.6665
$ws.Range("H67").Value = 11834
$ws.Range("I67").Value = 6446.6665
$ws.Range("K67").Value = 6446.6665
$ws.Range("M67").Value = -5588.6665
$ws.Range("H70").Value = 2284.2144
$ws.Range("I70").Value = 1730
$ws.Range("J70").Value = 2699.875
$ws.Range("K70").Value = 5190
$ws.Range("L70").Value = 8099.625
$ws.Range("M70").Value = -4920
$ws.Range("N70").Value = -8639.625
$ws.Range("H73").Value = 2284.2144
$ws.Range("I73").Value = 1730
$ws.Range("J73").Value = 2699.875
$ws.Range("K73").Value = 5190
$ws.Range("L73").Value = 8099.625
$ws.Range("M73").Value = -4254
$ws.Range("N73").Value = -9971.625
$ws.Range("H98").Value = 755.1429000000001
$ws.Range("I98").Value = 761.8125
$ws.Range("K98").Value = 761.8125
$ws.Range("M98").Value = 736.1875
$ws.Range("H122").Value = 755.1429000000001
$ws.Range("I122").Value = 761.8125
$ws.Range("K122").Value = 2285.4375
$ws.Range("M122").Value = 164.5625
$ws.Range("H138").Value = 3335.9119
$ws.Range("J138").Value = 3271.1
$ws.Range("L138").Value = 9813.299999999999
$ws.Range("N138").Value = -20093.3
$ws.Range("H141").Value = 7453.4165
$ws.Range("I141").Value = 7189.8335
$ws.Range("K141").Value = 21569.5005
$ws.Range("M141").Value = -16389.5005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4613.294
$ws.Range("I32").Value = 4405.28
$ws.Range("K32").Value = 4405.28
$ws.Range("M32").Value = -4118.28
$ws.Range("H110").Value = 6172.1
$ws.Range("I110").Value = 5230.2
$ws.Range("K110").Value = 5230.2
$ws.Range("M110").Value = -3185.2
$ws.Range("H132").Value = 1535.421
$ws.Range("I132").Value = 1555.6285
$ws.Range("K132").Value = 4666.8855
$ws.Range("M132").Value = -2136.8855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 69979
$ws.Range("J2").Value = 69979
$ws.Range("L2").Value = 69979
$ws.Range("N2").Value = -70205
$ws.Range("H22").Value = 958.3333
$ws.Range("I22").Value = 958.3333
$ws.Range("K22").Value = 958.3333
$ws.Range("M22").Value = -785.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 209.875
$ws.Range("I19").Value = 236.28572
$ws.Range("K19").Value = 236.28572
$ws.Range("M19").Value = -66.28572
$ws.Range("H24").Value = 209.875
$ws.Range("I24").Value = 236.28572
$ws.Range("K24").Value = 236.28572
$ws.Range("M24").Value = -66.28572
$ws.Range("H31").Value = 5010.5454
$ws.Range("I31").Value = 3465.25
$ws.Range("K31").Value = 3465.25
$ws.Range("M31").Value = -3170.25
$ws.Range("H34").Value = 5010.5454
$ws.Range("I34").Value = 3465.25
$ws.Range("K34").Value = 3465.25
$ws.Range("M34").Value = -3263.25
$ws.Range("H43").Value = 28500
$ws.Range("J43").Value = 28500
$ws.Range("L43").Value = 28500
$ws.Range("N43").Value = -28868
$ws.Range("H74").Value = 36496.668
$ws.Range("J74").Value = 36496.668
$ws.Range("L74").Value = 36496.668
$ws.Range("N74").Value = -38244.668
$ws.Range("H77").Value = 36496.668
$ws.Range("J77").Value = 36496.668
$ws.Range("L77").Value = 109490.004
$ws.Range("N77").Value = -118226.004
$ws.Range("H82").Value = 52999.668
$ws.Range("J82").Value = 52999.668
$ws.Range("L82").Value = 52999.668
$ws.Range("N82").Value = -53721.668
$ws.Range("H85").Value = 52999.668
$ws.Range("J85").Value = 52999.668
$ws.Range("L85").Value = 52999.668
$ws.Range("N85").Value = -55495.668
$ws.Range("H95").Value = 21431.5
$ws.Range("J95").Value = 21431.5
$ws.Range("L95").Value = 21431.5
$ws.Range("N95").Value = -26923.5
$ws.Range("H101").Value = 28500
$ws.Range("J101").Value = 28500
$ws.Range("L101").Value = 28500
$ws.Range("N101").Value = -34990

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2471.6667
$ws.Range("I132").Value = 1745
$ws.Range("K132").Value = 15705
$ws.Range("M132").Value = -13175

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 12402600
$ws.Range("I11").Value = 20171000
$ws.Range("K11").Value = 20171000
$ws.Range("M11").Value = -20170861
$ws.Range("H20").Value = 8585846
$ws.Range("J20").Value = 17184.6
$ws.Range("L20").Value = 17184.6
$ws.Range("N20").Value = -17674.6
$ws.Range("H101").Value = 62500
$ws.Range("J101").Value = 62500
$ws.Range("L101").Value = 62500
$ws.Range("N101").Value = -68990
$ws.Range("H126").Value = 3117.9048
$ws.Range("J126").Value = 4237.5
$ws.Range("L126").Value = 12712.5
$ws.Range("N126").Value = -17652.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 1000025
$ws.Range("I19").Value = 50
$ws.Range("K19").Value = 50
$ws.Range("M19").Value = 120
$ws.Range("H40").Value = 3562.6086
$ws.Range("I40").Value = 3055.2942
$ws.Range("K40").Value = 3055.2942
$ws.Range("M40").Value = -2919.2942
$ws.Range("H93").Value = 18614.5
$ws.Range("I93").Value = 922.25
$ws.Range("K93").Value = 922.25
$ws.Range("M93").Value = 325.75
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H132").Value = 7064.3184
$ws.Range("I132").Value = 8672
$ws.Range("J132").Value = 3619.2856
$ws.Range("K132").Value = 26016
$ws.Range("L132").Value = 10857.8568
$ws.Range("M132").Value = -23486
$ws.Range("N132").Value = -15917.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 6982
$ws.Range("I15").Value = 6982
$ws.Range("K15").Value = 6982
$ws.Range("M15").Value = -6694
$ws.Range("H100").Value = 1392.8
$ws.Range("I100").Value = 1354.6666
$ws.Range("K100").Value = 2709.3332
$ws.Range("M100").Value = -2168.3332
$ws.Range("H107").Value = 1135.1818
$ws.Range("J107").Value = 1267.6
$ws.Range("L107").Value = 3802.8
$ws.Range("N107").Value = -7642.799999999999
$ws.Range("H126").Value = 7424
$ws.Range("I126").Value = 4033.1667
$ws.Range("J126").Value = 10330.429
$ws.Range("K126").Value = 12099.5001
$ws.Range("L126").Value = 30991.287
$ws.Range("M126").Value = -9629.500100000001
$ws.Range("N126").Value = -35931.287
$ws.Range("H136").Value = 2869.3333
$ws.Range("I136").Value = 2762.7837
$ws.Range("K136").Value = 8288.3511
$ws.Range("M136").Value = -5738.3511
